$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("all_tools")
$ws.Range("L5").Value = 0.3140610795011266
$ws.Range("G9").Value = 1127
$ws.Range("I9").Value = -0.1624207945913728
$ws.Range("J9").Value = 0.01995735599907151
$ws.Range("K9").Value = -0.2433748586082731
$ws.Range("L9").Value = 0.01468902726730107
$ws.Range("G25").Value = 42
$ws.Range("I25").Value = -0.2089183206099184
$ws.Range("J25").Value = 0.2858071051160878
$ws.Range("K25").Value = -0.2761088800464458
$ws.Range("L25").Value = 0.3006008438154604
$ws.Range("G26").Value = 42
$ws.Range("I26").Value = -0.1907515101220994
$ws.Range("J26").Value = 0.3297832606788835
$ws.Range("K26").Value = -0.2791264634349316
$ws.Range("L26").Value = 0.2951316609645892
$ws.Range("G27").Value = 42
$ws.Range("I27").Value = 0.3648596387010686
$ws.Range("J27").Value = 0.06317307481323101
$ws.Range("K27").Value = 0.4446662347157511
$ws.Range("L27").Value = 0.08439298705280707
$ws.Range("G28").Value = 42
$ws.Range("I28").Value = -0.1563850753205831
$ws.Range("J28").Value = 0.4285752941987943
$ws.Range("K28").Value = -0.2133681077431026
$ws.Range("L28").Value = 0.4275202599269713
$ws.Range("G29").Value = 42
$ws.Range("I29").Value = -0.1725846996342804
$ws.Range("J29").Value = 0.377919481405675
$ws.Range("K29").Value = -0.1931253368630878
$ws.Range("L29").Value = 0.4736035833560223

$ws = $wb.Worksheets.Item("checker_framework")
$ws.Range("F9").Value = 17
$ws.Range("G9").Value = 50
$ws.Range("I9").Value = -0.2260489959509954
$ws.Range("J9").Value = 0.004741931174898608
$ws.Range("K9").Value = -0.2832375756589699
$ws.Range("L9").Value = 0.004297109542258244
$ws.Range("L11").Value = 0.004913052567611288

$ws = $wb.Worksheets.Item("typestate_checker")
$ws.Range("G9").Value = 326
$ws.Range("I9").Value = -0.1509654696395568
$ws.Range("J9").Value = 0.03611918064522761
$ws.Range("K9").Value = -0.2132165514106953
$ws.Range("L9").Value = 0.0331785702669911
$ws.Range("L14").Value = 0.4230203924441358

$ws = $wb.Worksheets.Item("infer")
$ws.Range("L12").Value = 0.03867934687031339

$ws = $wb.Worksheets.Item("openjml")
$ws.Range("L7").Value = 0.06043495620092659
$ws.Range("L9").Value = 0.09653017580355105
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 32
$ws.Range("I25").Value = -0.3244079990200284
$ws.Range("J25").Value = 0.1020809611324846
$ws.Range("K25").Value = -0.4824165900576836
$ws.Range("L25").Value = 0.05842088299110668
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 32
$ws.Range("I26").Value = -0.4170959987400364
$ws.Range("J26").Value = 0.03555790569992093
$ws.Range("K26").Value = -0.5355127556300702
$ws.Range("L26").Value = 0.03253137976639191
$ws.Range("F27").Value = 14
$ws.Range("G27").Value = 32
$ws.Range("I27").Value = 0.2047685867790963
$ws.Range("J27").Value = 0.3037024766857317
$ws.Range("K27").Value = 0.2170954499333818
$ws.Range("L27").Value = 0.4192967698838219
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 32
$ws.Range("I28").Value = 0.08448190755542286
$ws.Range("J28").Value = 0.6731676659685988
$ws.Range("K28").Value = 0.1201996807657877
$ws.Range("L28").Value = 0.6574684200274921
$ws.Range("F29").Value = 14
$ws.Range("G29").Value = 32
$ws.Range("I29").Value = 0.009268799972000809
$ws.Range("J29").Value = 0.9627444106905073
$ws.Range("K29").Value = 0.01365329971861369
$ws.Range("L29").Value = 0.9599750222894384

# Column J (10) width increased by one character width
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(10).ColumnWidth + 1
